# Updated cryptos list with GitHub Actions run data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value (applied as literal text, matching the
# original "Price"/"Volume(1h)" columns which store formatted strings, not numbers)
$updates = @{
    'D2' = '61.501.86'
    'E2' = '  +0.58%  '
    'D3' = '3.445.07'
    'E3' = '  +1.41%  '
    'E4' = '  -0.04%  '
    'D5' = '577.71'
    'E5' = '  +0.92%  '
    'D6' = '145.08'
    'E6' = '  +4.75%  '
    'D7' = '3.444.63'
    'E7' = '  +1.42%  '
    'E8' = '  +0.05%  '
    'E9' = '  +2.33%  '
    'E10' = '  -0.17%  '
    'E11' = '  +3.69%  '
    'E12' = '  +2.48%  '
    'D13' = '4.032.62'
    'E13' = '  +1.40%  '
    'D14' = '28.35'
    'E14' = '  +6.28%  '
    'E15' = '  -0.45%  '
    'D16' = '0.0000173'
    'E16' = '  +1.26%  '
    'D17' = '3.456.45'
    'E17' = '  +1.64%  '
    'D18' = '61.651.60'
    'E18' = '  +0.69%  '
    'E19' = '  +6.82%  '
    'D20' = '14.34'
    'E20' = '  +3.60%  '
    'D21' = '9.42'
    'E21' = '  +0.98%  '
    'D22' = '402.13'
    'E22' = '  +7.07%  '
    'E23' = '  +3.17%  '
    'D24' = '74.48'
    'E24' = '  +4.68%  '
    'E25' = '  +0.13%  '
    'E26' = '  +0.71%  '
    'E27' = '  +0.32%  '
    'D28' = '3.589.18'
    'D29' = '0.182'
    'E29' = '  +4.46%  '
    'E30' = '  +3.30%  '
    'E32' = '  +2.14%  '
    'E33' = '  +2.11%  '
    'E34' = '  -9.47%  '
    'E36' = '  +2.06%  '
    'B37' = 'RenzoRestakedETH'
    'C37' = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
    'D37' = '3.472.84'
    'E37' = '  +1.64%  '
    'B38' = 'Aptos'
    'C38' = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
    'D38' = '7.04'
    'E38' = '  +2.80%  '
    'E39' = '  +0.31%  '
    'E40' = '  +0.72%  '
    'D41' = '167.07'
    'E41' = '  +0.67%  '
    'D42' = '0.0790'
    'E42' = '  +2.87%  '
    'D43' = '27.19'
    'E43' = '  +4.15%  '
    'D44' = '0.802'
    'E44' = '  +3.23%  '
    'D45' = '4.52'
    'E45' = '  +3.20%  '
    'B46' = 'FirstDigitalUSD'
    'C46' = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
    'D46' = '1.00'
    'E46' = '  -0.02%  '
    'B47' = 'Stacks'
    'C47' = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
    'D47' = '1.73'
    'E47' = '  -1.09%  '
    'D48' = '42.37'
    'E48' = '  +1.14%  '
    'D49' = '2.613.00'
    'E49' = '  +3.67%  '
    'E50' = '  -2.01%  '
    'D51' = '6.95'
    'E51' = '  +2.57%  '
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    # Force text format so numeric-looking strings (e.g. "1.00", "577.71")
    # are kept as text instead of being parsed into numbers.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    # Restore the default "Normal" style so no stray cell style/number format
    # is left behind on cells that originally had no explicit style.
    $cell.Style = "Normal"
}
